# Add a new FAV_NUMBER column (F) with per-row values, matching the
# "Write field values using proper java reflection" commit which adds a
# sample int-typed field/column to the test fixture.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "FAV_NUMBER"
$ws.Range("F2").Value = 3
$ws.Range("F3").Value = 2
$ws.Range("F4").Value = 4
$ws.Range("F5").Value = 10
$ws.Range("F6").Value = 1

# Move the active selection to match the new last-written cell.
[void]$ws.Range("F4").Select()
